# LF Energy High Level Overview Deck - February 2025 update
# - Re-style the 6 "member benefits" comparison tables with the new
#   built-in table style.
# - Swap the two theme colour palettes ("Simple Light" <-> "Geometric").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Tables: switch every comparison-table on slides 9, 13-16 and 18
#    from the old custom style to the new built-in style.
# ---------------------------------------------------------------------
$newTableStyle = "{591A0ED6-D4D7-4D79-807F-AB3D0AB6F1A9}"

$tableShapes = @{
    9  = 3
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    18 = 3
}

foreach ($slideIdx in $tableShapes.Keys) {
    $shapeIdx = $tableShapes[$slideIdx]
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item($shapeIdx)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle, $true)
    }
}

# ---------------------------------------------------------------------
# 2) Theme: the deck's colour scheme is being swapped from the old
#    "Simple Light" palette to the "Geometric" (LF Energy 2023) palette.
# ---------------------------------------------------------------------
function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $c = $scheme.Colors($index)
    $c.RGB = $r + ($g * 256) + ($b * 65536)
}

$scheme = $p.Slides.Item(1).ThemeColorScheme

# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
Set-ThemeColor $scheme 1  0x22 0x22 0x22   # dk1
Set-ThemeColor $scheme 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $scheme 3  0x43 0x43 0x43   # dk2
Set-ThemeColor $scheme 4  0x99 0x99 0x99   # lt2
Set-ThemeColor $scheme 5  0x00 0x37 0x78   # accent1
Set-ThemeColor $scheme 6  0x00 0x94 0xFF   # accent2
Set-ThemeColor $scheme 7  0x5B 0x1D 0xE7   # accent3
Set-ThemeColor $scheme 8  0x12 0xE2 0xE2   # accent4
Set-ThemeColor $scheme 9  0xFF 0x00 0xAA   # accent5
Set-ThemeColor $scheme 10 0xAC 0xDE 0x1F   # accent6
Set-ThemeColor $scheme 11 0x00 0x77 0xCC   # hlink
Set-ThemeColor $scheme 12 0xF0 0x62 0x92   # folHlink
